$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = 0.8950246002264373
$ws.Cells.Item(2, 10).Value = 0.8950246002264374
$ws.Cells.Item(2, 13).Value = 20.79761
$ws.Cells.Item(2, 14).Value = 41.59522
$ws.Cells.Item(2, 15).Value = 0.1221346813081061
$ws.Cells.Item(2, 16).Value = 0.08652867365760288
$ws.Cells.Item(2, 17).Value = 35.10719065186333
$ws.Cells.Item(2, 18).Value = 210.64314391118
$ws.Cells.Item(2, 19).Value = 0.109313544311571
$ws.Cells.Item(2, 20).Value = 0.07744529154851988
$ws.Cells.Item(3, 9).Value = 0.8950246002264373
$ws.Cells.Item(3, 10).Value = 0.8950246002264374
$ws.Cells.Item(3, 13).Value = 13.56641266666666
$ws.Cells.Item(3, 14).Value = 40.69923799999999
$ws.Cells.Item(3, 15).Value = 0.0796692258166966
$ws.Cells.Item(3, 16).Value = 0.08466480242237233
$ws.Cells.Item(3, 19).Value = 0.07130591698693863
$ws.Cells.Item(3, 20).Value = 0.0757770809413341
$ws.Cells.Item(4, 9).Value = 0.8950246002264373
$ws.Cells.Item(4, 10).Value = 0.8950246002264374
$ws.Cells.Item(4, 13).Value = 49.80229833333333
$ws.Cells.Item(4, 14).Value = 149.406895
$ws.Cells.Item(4, 15).Value = 0.2924657129041698
$ws.Cells.Item(4, 16).Value = 0.3108044736787241
$ws.Cells.Item(4, 17).Value = 84.06825507783388
$ws.Cells.Item(4, 18).Value = 756.614295700505
$ws.Cells.Item(4, 19).Value = 0.2617640077719946
$ws.Cells.Item(4, 20).Value = 0.2781776498028883
$ws.Cells.Item(5, 9).Value = 0.8950246002264373
$ws.Cells.Item(5, 10).Value = 0.8950246002264374
$ws.Cells.Item(5, 13).Value = 9.34483
$ws.Cells.Item(5, 14).Value = 18.68966
$ws.Cells.Item(5, 15).Value = 0.05487783615177078
$ws.Cells.Item(5, 16).Value = 0.03887926283143963
$ws.Cells.Item(5, 17).Value = 15.77444371825667
$ws.Cells.Item(5, 18).Value = 94.64666230953999
$ws.Cells.Item(5, 19).Value = 0.04911701336303056
$ws.Cells.Item(5, 20).Value = 0.03479789667280784
$ws.Cells.Item(6, 9).Value = 0.8950246002264373
$ws.Cells.Item(6, 10).Value = 0.8950246002264374
$ws.Cells.Item(6, 13).Value = 37.76134866666666
$ws.Cells.Item(6, 14).Value = 113.284046
$ws.Cells.Item(6, 15).Value = 0.2217548211149075
$ws.Cells.Item(6, 16).Value = 0.2356597283761661
$ws.Cells.Item(6, 17).Value = 63.74265441616377
$ws.Cells.Item(6, 18).Value = 573.6838897454739
$ws.Cells.Item(6, 19).Value = 0.1984760201166552
$ws.Cells.Item(6, 20).Value = 0.2109212541793489
$ws.Cells.Item(7, 9).Value = 0.8950246002264373
$ws.Cells.Item(7, 10).Value = 0.8950246002264374
$ws.Cells.Item(7, 13).Value = 39.011729
$ws.Cells.Item(7, 14).Value = 117.035187
$ws.Cells.Item(7, 15).Value = 0.2290977227043493
$ws.Cells.Item(7, 16).Value = 0.2434630590336949
$ws.Cells.Item(7, 17).Value = 65.85334601725033
$ws.Cells.Item(7, 18).Value = 592.680114155253
$ws.Cells.Item(7, 19).Value = 0.2050480976762474
$ws.Cells.Item(7, 20).Value = 0.2179054270815383
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.175453
$ws.Cells.Item(8, 8).Value = 0.526359
$ws.Cells.Item(8, 9).Value = 0.0930278797853264
$ws.Cells.Item(8, 10).Value = 0.09302787978532641
$ws.Cells.Item(8, 13).Value = 20.79761
$ws.Cells.Item(8, 14).Value = 41.59522
$ws.Cells.Item(8, 15).Value = 0.1221346813081061
$ws.Cells.Item(8, 16).Value = 0.08652867365760288
$ws.Cells.Item(8, 17).Value = 3.64900306733
$ws.Cells.Item(8, 18).Value = 21.89401840398
$ws.Cells.Item(8, 19).Value = 0.01136193045034965
$ws.Cells.Item(8, 20).Value = 0.00804957905100322
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.175453
$ws.Cells.Item(9, 8).Value = 0.526359
$ws.Cells.Item(9, 9).Value = 0.0930278797853264
$ws.Cells.Item(9, 10).Value = 0.09302787978532641
$ws.Cells.Item(9, 13).Value = 13.56641266666666
$ws.Cells.Item(9, 14).Value = 40.69923799999999
$ws.Cells.Item(9, 15).Value = 0.0796692258166966
$ws.Cells.Item(9, 16).Value = 0.08466480242237233
$ws.Cells.Item(9, 17).Value = 2.380267801604666
$ws.Cells.Item(9, 18).Value = 21.422410214442
$ws.Cells.Item(9, 19).Value = 0.007411459161865673
$ws.Cells.Item(9, 20).Value = 0.007876187061796866
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.175453
$ws.Cells.Item(10, 8).Value = 0.526359
$ws.Cells.Item(10, 9).Value = 0.0930278797853264
$ws.Cells.Item(10, 10).Value = 0.09302787978532641
$ws.Cells.Item(10, 13).Value = 49.80229833333333
$ws.Cells.Item(10, 14).Value = 149.406895
$ws.Cells.Item(10, 15).Value = 0.2924657129041698
$ws.Cells.Item(10, 16).Value = 0.3108044736787241
$ws.Cells.Item(10, 17).Value = 8.737962649478334
$ws.Cells.Item(10, 18).Value = 78.641663845305
$ws.Cells.Item(10, 19).Value = 0.02720746518137889
$ws.Cells.Item(10, 20).Value = 0.02891348121412599
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.175453
$ws.Cells.Item(11, 8).Value = 0.526359
$ws.Cells.Item(11, 9).Value = 0.0930278797853264
$ws.Cells.Item(11, 10).Value = 0.09302787978532641
$ws.Cells.Item(11, 13).Value = 9.34483
$ws.Cells.Item(11, 14).Value = 18.68966
$ws.Cells.Item(11, 15).Value = 0.05487783615177078
$ws.Cells.Item(11, 16).Value = 0.03887926283143963
$ws.Cells.Item(11, 17).Value = 1.63957845799
$ws.Cells.Item(11, 18).Value = 9.837470747940001
$ws.Cells.Item(11, 19).Value = 0.005105168744405771
$ws.Cells.Item(11, 20).Value = 0.003616855388825275
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.175453
$ws.Cells.Item(12, 8).Value = 0.526359
$ws.Cells.Item(12, 9).Value = 0.0930278797853264
$ws.Cells.Item(12, 10).Value = 0.09302787978532641
$ws.Cells.Item(12, 13).Value = 37.76134866666666
$ws.Cells.Item(12, 14).Value = 113.284046
$ws.Cells.Item(12, 15).Value = 0.2217548211149075
$ws.Cells.Item(12, 16).Value = 0.2356597283761661
$ws.Cells.Item(12, 17).Value = 6.625341907612666
$ws.Cells.Item(12, 18).Value = 59.628077168514
$ws.Cells.Item(12, 19).Value = 0.02062938084049417
$ws.Cells.Item(12, 20).Value = 0.02192292488162065
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.175453
$ws.Cells.Item(13, 8).Value = 0.526359
$ws.Cells.Item(13, 9).Value = 0.0930278797853264
$ws.Cells.Item(13, 10).Value = 0.09302787978532641
$ws.Cells.Item(13, 13).Value = 39.011729
$ws.Cells.Item(13, 14).Value = 117.035187
$ws.Cells.Item(13, 15).Value = 0.2290977227043493
$ws.Cells.Item(13, 16).Value = 0.2434630590336949
$ws.Cells.Item(13, 17).Value = 6.844724888237001
$ws.Cells.Item(13, 18).Value = 61.60252399413301
$ws.Cells.Item(13, 19).Value = 0.02131247540683225
$ws.Cells.Item(13, 20).Value = 0.0226488521879544
$ws.Cells.Item(14, 7).Value = 0.02253333333333333
$ws.Cells.Item(14, 8).Value = 0.06759999999999999
$ws.Cells.Item(14, 9).Value = 0.01194751998823629
$ws.Cells.Item(14, 10).Value = 0.01194751998823629
$ws.Cells.Item(14, 13).Value = 20.79761
$ws.Cells.Item(14, 14).Value = 41.59522
$ws.Cells.Item(14, 15).Value = 0.1221346813081061
$ws.Cells.Item(14, 16).Value = 0.08652867365760288
$ws.Cells.Item(14, 17).Value = 0.4686394786666666
$ws.Cells.Item(14, 18).Value = 2.811836872
$ws.Cells.Item(14, 19).Value = 0.001459206546185467
$ws.Cells.Item(14, 20).Value = 0.001033803058079785
$ws.Cells.Item(15, 7).Value = 0.02253333333333333
$ws.Cells.Item(15, 8).Value = 0.06759999999999999
$ws.Cells.Item(15, 9).Value = 0.01194751998823629
$ws.Cells.Item(15, 10).Value = 0.01194751998823629
$ws.Cells.Item(15, 13).Value = 13.56641266666666
$ws.Cells.Item(15, 14).Value = 40.69923799999999
$ws.Cells.Item(15, 15).Value = 0.0796692258166966
$ws.Cells.Item(15, 16).Value = 0.08466480242237233
$ws.Cells.Item(15, 17).Value = 0.3056964987555555
$ws.Cells.Item(15, 18).Value = 2.751268488799999
$ws.Cells.Item(15, 19).Value = 0.0009518496678922931
$ws.Cells.Item(15, 20).Value = 0.001011534419241369
$ws.Cells.Item(16, 7).Value = 0.02253333333333333
$ws.Cells.Item(16, 8).Value = 0.06759999999999999
$ws.Cells.Item(16, 9).Value = 0.01194751998823629
$ws.Cells.Item(16, 10).Value = 0.01194751998823629
$ws.Cells.Item(16, 13).Value = 49.80229833333333
$ws.Cells.Item(16, 14).Value = 149.406895
$ws.Cells.Item(16, 15).Value = 0.2924657129041698
$ws.Cells.Item(16, 16).Value = 0.3108044736787241
$ws.Cells.Item(16, 17).Value = 1.122211789111111
$ws.Cells.Item(16, 18).Value = 10.099906102
$ws.Cells.Item(16, 19).Value = 0.003494239950796345
$ws.Cells.Item(16, 20).Value = 0.003713342661709816
$ws.Cells.Item(17, 7).Value = 0.02253333333333333
$ws.Cells.Item(17, 8).Value = 0.06759999999999999
$ws.Cells.Item(17, 9).Value = 0.01194751998823629
$ws.Cells.Item(17, 10).Value = 0.01194751998823629
$ws.Cells.Item(17, 13).Value = 9.34483
$ws.Cells.Item(17, 14).Value = 18.68966
$ws.Cells.Item(17, 15).Value = 0.05487783615177078
$ws.Cells.Item(17, 16).Value = 0.03887926283143963
$ws.Cells.Item(17, 17).Value = 0.2105701693333333
$ws.Cells.Item(17, 18).Value = 1.263421016
$ws.Cells.Item(17, 19).Value = 0.0006556540443344374
$ws.Cells.Item(17, 20).Value = 0.0004645107698065172
$ws.Cells.Item(18, 7).Value = 0.02253333333333333
$ws.Cells.Item(18, 8).Value = 0.06759999999999999
$ws.Cells.Item(18, 9).Value = 0.01194751998823629
$ws.Cells.Item(18, 10).Value = 0.01194751998823629
$ws.Cells.Item(18, 13).Value = 37.76134866666666
$ws.Cells.Item(18, 14).Value = 113.284046
$ws.Cells.Item(18, 15).Value = 0.2217548211149075
$ws.Cells.Item(18, 16).Value = 0.2356597283761661
$ws.Cells.Item(18, 17).Value = 0.8508890566222221
$ws.Cells.Item(18, 18).Value = 7.658001509599998
$ws.Cells.Item(18, 19).Value = 0.00264942015775812
$ws.Cells.Item(18, 20).Value = 0.002815549315196579
$ws.Cells.Item(19, 7).Value = 0.02253333333333333
$ws.Cells.Item(19, 8).Value = 0.06759999999999999
$ws.Cells.Item(19, 9).Value = 0.01194751998823629
$ws.Cells.Item(19, 10).Value = 0.01194751998823629
$ws.Cells.Item(19, 13).Value = 39.011729
$ws.Cells.Item(19, 14).Value = 117.035187
$ws.Cells.Item(19, 15).Value = 0.2290977227043493
$ws.Cells.Item(19, 16).Value = 0.2434630590336949
$ws.Cells.Item(19, 17).Value = 0.8790642934666667
$ws.Cells.Item(19, 18).Value = 7.911578641199999
$ws.Cells.Item(19, 19).Value = 0.002737149621269628
$ws.Cells.Item(19, 20).Value = 0.002908779764202221
